# Apply new values to column A (rows 1-99) of the active sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    25,
    11,
    29,
    3,
    11,
    4,
    27,
    28,
    6,
    11,
    12,
    22,
    13,
    15,
    12,
    7,
    20,
    16,
    10,
    12,
    2,
    5,
    17,
    26,
    28,
    25,
    13,
    13,
    15,
    10,
    16,
    2,
    30,
    12,
    29,
    20,
    12,
    24,
    7,
    6,
    3,
    7,
    14,
    28,
    14,
    19,
    18,
    2,
    20,
    20,
    27,
    11,
    25,
    21,
    4,
    13,
    3,
    16,
    25,
    22,
    20,
    21,
    21,
    28,
    29,
    1,
    21,
    21,
    25,
    4,
    20,
    17,
    24,
    30,
    23,
    17,
    27,
    4,
    7,
    17,
    22,
    23,
    10,
    5,
    14,
    29,
    1,
    21,
    8,
    4,
    26,
    18,
    8,
    18,
    29,
    21,
    17,
    9,
    25
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
